# Replace "null" (empty) I/J (lat/long) cells with 1000 so they can be
# filtered out later.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1")

$rows = @(2, 3, 9, 10, 11, 22, 23, 27, 34, 37, 38, 42, 45, 46, 65, 70, 71, 83, 84, 87, 112, 139)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = 1000
    $ws.Cells.Item($r, 10).Value = 1000
}

# Keep the active selection roughly where the author left off editing.
$ws.Range("I158").Select()
